# Generate Report for Handoff
# The 88703ad9-6223-4c18-89e1-d3203bc19335.md file moved from
# "Handed back: in sync with en-US" to "Ready for handoff", and the
# "Latest Handoff Datetime" columns were refreshed with the new handoff run.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet: status for the 88703ad9 file (row 3) in both locales ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-03-03 11:04:35"
$wsZhCn.Range("B3").Value = $newStatus
$wsZhCn.Range("D3").Value = "2016-03-03 11:04:35"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-03-03 11:04:47"
$wsDeDe.Range("B3").Value = $newStatus
$wsDeDe.Range("D3").Value = "2016-03-03 11:04:47"
